$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 12-24 (their content is no longer needed)
$ws.Range("A12:A24").ClearContents()

# Write the new values into A1:A11
$ws.Range("A1").Value = "above "
$ws.Range("A2").Value = "above "
$ws.Range("A3").Value = "apple"
$ws.Range("A4").Value = "banana"
$ws.Range("A5").Value = "coconut"
$ws.Range("A6").Value = "elephant"
$ws.Range("A7").Value = "frozen"
$ws.Range("A8").Value = "fancy "
$ws.Range("A9").Value = "hello "
$ws.Range("A10").Value = "yummy"
$ws.Range("A11").Value = "zombie"

# A2 gets a "vertical = top" alignment style
$ws.Range("A2").VerticalAlignment = -4160

# Update the selection to A2
$ws.Range("A2").Select()

# Set page orientation to Portrait
$ws.PageSetup.Orientation = 1
